$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet lists cryptocurrencies with their Price (column D) and
# Volume(1h) change (column E) as plain text cells. This updates those two
# columns with refreshed values for the rows that changed.
#
# Column D values like "594.07" look numeric, so Excel would normally
# auto-convert them into Number cells on assignment. To keep them as text
# (matching the original file), each Price cell's NumberFormat is
# temporarily set to Text ("@") before the value is written, and then
# ClearFormats() restores the cell's original (default) formatting
# afterwards so no stray number format is left behind.
#
# Column E values already start/end with extra spaces, which is enough for
# Excel to treat them as text, so no special handling is needed there.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "68.141.58"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +2.10%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.534.03"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.32%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "594.07"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +1.52%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "177.18"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  -0.03%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.532"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +0.91%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.533.18"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("E11").Value = "  +2.58%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "5.13"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  -1.12%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "26.84"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.17%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.992.61"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +1.85%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "68.185.35"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +2.37%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.538.26"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.67%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "8.00"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.63%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "11.57"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +2.02%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "366.80"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  -2.47%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "70.93"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("E28").Value = "  -0.03%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.655.38"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  +1.47%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "542.01"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +2.61%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "8.29"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.63%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "157.25"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("E39").Value = "  +1.09%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.356"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("E43").Value = "  +1.91%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  -0.06%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "147.43"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.48%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.561"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E49").Value = "  +2.43%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.71"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  -0.15%  "
